# Apply the price/volume data refresh captured in the commit
# "Updated symbol list on Fri Jan  6 21:40:11 UTC 2023 with GitHub Actions".
# The Price (column D) and Volume(1h) (column E) text values are refreshed
# in-place for the affected rows. A leading apostrophe is used so Excel
# stores the numeric/percentage-looking strings as literal text, matching
# the original inline-string cell contents instead of converting them to
# numbers or percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'259.29"
$ws.Range("E2").Value = "'0.78%"
$ws.Range("D3").Value = "'27.01"
$ws.Range("E3").Value = "'-0.18%"
$ws.Range("D4").Value = "'4.698"
$ws.Range("E4").Value = "'0.12%"
$ws.Range("D5").Value = "'0.06032"
$ws.Range("E5").Value = "'2.57%"
$ws.Range("D6").Value = "'6.672"
$ws.Range("E6").Value = "'0.40%"
$ws.Range("D7").Value = "'0.8589"
$ws.Range("E7").Value = "'0.12%"
$ws.Range("D8").Value = "'0.9228"
$ws.Range("E8").Value = "'-3.94%"
$ws.Range("E9").Value = "'-0.86%"
$ws.Range("D10").Value = "'0.04952"
$ws.Range("E10").Value = "'25.25%"
$ws.Range("D11").Value = "'0.07088"
$ws.Range("E11").Value = "'-0.11%"
$ws.Range("D12").Value = "'0.03088"
$ws.Range("E12").Value = "'-2.97%"
$ws.Range("D13").Value = "'0.09126"
$ws.Range("E13").Value = "'-0.50%"
$ws.Range("D14").Value = "'0.001529"
$ws.Range("E14").Value = "'-0.63%"
$ws.Range("D15").Value = "'0.0006039"
$ws.Range("E15").Value = "'-0.40%"
$ws.Range("D16").Value = "'0.006081"
$ws.Range("E16").Value = "'-2.02%"
$ws.Range("D17").Value = "'3.465"
$ws.Range("E17").Value = "'-1.49%"
$ws.Range("D18").Value = "'3.168"
$ws.Range("D19").Value = "'2.166"
$ws.Range("E19").Value = "'-1.83%"
$ws.Range("E20").Value = "'0.43%"
$ws.Range("D21").Value = "'0.1297"
$ws.Range("D22").Value = "'4.117"
$ws.Range("E22").Value = "'6.84%"
$ws.Range("D23").Value = "'0.04238"
$ws.Range("E23").Value = "'0.30%"
$ws.Range("E24").Value = "'-0.54%"
$ws.Range("D25").Value = "'0.004038"
$ws.Range("E26").Value = "'-0.11%"
$ws.Range("E27").Value = "'-21.43%"
$ws.Range("D40").Value = "'0.03846"
$ws.Range("E40").Value = "'0.43%"
$ws.Range("E41").Value = "'1.22%"
$ws.Range("E42").Value = "'-35.57%"
$ws.Range("D43").Value = "'0.01514"
$ws.Range("E43").Value = "'32.36%"
$ws.Range("E44").Value = "'15.67%"
$ws.Range("D45").Value = "'0.00005092"
$ws.Range("E45").Value = "'-6.74%"
$ws.Range("E46").Value = "'-0.10%"
$ws.Range("D47").Value = "'0.05452"
$ws.Range("E47").Value = "'-9.18%"
$ws.Range("E48").Value = "'-11.90%"
$ws.Range("E49").Value = "'-0.10%"
$ws.Range("E50").Value = "'-0.10%"
